$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 7, pushing current rows 7-10 down to 9-12.
$ws.Rows.Item(7).Resize(2).Insert()

# New row 7 - Primera, updated 2021-12-06 data
$ws.Cells.Item(7, 1).Value = 3
$ws.Cells.Item(7, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44536
$ws.Cells.Item(7, 5).Value = 5
$ws.Cells.Item(7, 6).Value = 100112043
$ws.Cells.Item(7, 7).Value = "Pepino dulce"
$ws.Cells.Item(7, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 87
$ws.Cells.Item(7, 11).Value = 22000
$ws.Cells.Item(7, 12).Value = 22000
$ws.Cells.Item(7, 13).Value = 22000
$ws.Cells.Item(7, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 1222
$ws.Cells.Item(7, 17).Value = 18
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# New row 8 - Segunda, updated 2021-12-06 data
$ws.Cells.Item(8, 1).Value = 3
$ws.Cells.Item(8, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 44536
$ws.Cells.Item(8, 5).Value = 5
$ws.Cells.Item(8, 6).Value = 100112043
$ws.Cells.Item(8, 7).Value = "Pepino dulce"
$ws.Cells.Item(8, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(8, 9).Value = "Segunda"
$ws.Cells.Item(8, 10).Value = 80
$ws.Cells.Item(8, 11).Value = 20000
$ws.Cells.Item(8, 12).Value = 20000
$ws.Cells.Item(8, 13).Value = 20000
$ws.Cells.Item(8, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 1111
$ws.Cells.Item(8, 17).Value = 18
$ws.Cells.Item(8, 18).Value = "Hortaliza"
